$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price/volume figures and row reordering per latest data pull.

$ws.Range("D2").Value = "42.849.23"
$ws.Range("E2").Value = "  +1.09%  "

$ws.Range("D3").Value = "2.288.18"
$ws.Range("E3").Value = "  -0.58%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.01"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.51%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.79"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.43%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "103.86"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.18%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.623"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.33%  "

$ws.Range("E8").Value = "  +0.36%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.600"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.65%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.26"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.20%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0902"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.84%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.41"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.53%  "

$ws.Range("E13").Value = "  +2.29%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.00"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.52%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.21"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.76%  "

$ws.Range("D16").Value = "2.642.39"
$ws.Range("E16").Value = "  -0.25%  "

$ws.Range("D17").Value = "2.282.57"
$ws.Range("E17").Value = "  -1.13%  "

$ws.Range("D18").Value = "42.783.74"
$ws.Range("E18").Value = "  +0.99%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.42"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.81%  "

$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000105"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.74%  "

$ws.Range("B21").Value = "InternetComputer(DFINITY)"
$ws.Range("C21").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.40"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +22.51%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.87"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.69%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.55"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.16%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "261.99"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.20%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.20"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.69%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.01"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.38%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.86"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.17%  "

$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.35"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.06%  "

$ws.Range("B29").Value = "Filecoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.04"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +19.63%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.25"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.28%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "37.37"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.28%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "166.53"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.92%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0870"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.27%  "

$ws.Range("E34").Value = "  -4.17%  "

$ws.Range("E35").Value = "  -0.60%  "

$ws.Range("E36").Value = "  -2.66%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.54"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.16%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.86"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.50%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0349"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.73%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.66"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.68%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.56"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.51%  "

$ws.Range("E42").Value = "  +1.10%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "69.18"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.25%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.01"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.45%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "91.73"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.05%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.13"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.78%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "113.65"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.85%  "

$ws.Range("D48").Value = "1.724.73"
$ws.Range("E48").Value = "  +7.72%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "78.76"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.01%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.76"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.69%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.20"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.64%  "

